$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column G (shifts FOTO/L-P/ID MATERI columns right)
$ws.Columns("G").EntireColumn.Insert()

# New header text for the inserted column
$ws.Range("G1").Value = "NO BILLKEY"

# Clear the inherited formatting/content in the new column's data rows (2-4) -
# they should remain empty, matching the sibling empty columns (D/E)
$ws.Range("G2:G4").Clear()

# Match column width to other wide text columns (NAMA/TMP LAHIR column width)
$ws.Columns("G").ColumnWidth = 35.86

# Update selection to reflect where the user was working
$ws.Range("G10").Select()

Write-Output "done"
